$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row for "2022-Q3" data.
#    The existing row 2 ("2021-Q2") is pushed down to row 3, and the
#    previous row 3 ("2021-Q1") is pushed down to row 4. Row 2 is then
#    overwritten with the new "2022-Q3" figures.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Make room for a new row at position 3 (shifts old row3 "2021-Q1" -> row4)
$wsTotal.Rows.Item(3).Insert()

# Re-apply the index-column style to the newly inserted A3 cell (Insert()
# otherwise invents a blended border/format that does not match A2/A4).
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# Move the old row 2 ("2021-Q2": 1, 0) down into the new row 3.
$wsTotal.Range("A2:D2").Copy()
$wsTotal.Range("A3:D3").PasteSpecial(-4104)
$wsTotal.Application.CutCopyMode = $false

# Fix up the running index in column A for the shifted rows.
$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(4,1).Value = 2

# Overwrite row 2 with the brand-new "2022-Q3" summary figures.
$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 3
$wsTotal.Cells.Item(2,4).Value = 0.08

# ---------------------------------------------------------------------------
# 2. Add a brand-new worksheet "2022-Q3" positioned right after "总计" and
#    before "2021-Q2" (mirrors the existing per-quarter fund-holding sheets).
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item(2)
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Keep the fund-code / ratio columns as text (so leading zeros such as in
# "006165" and values like "0.0530" are preserved verbatim), same as the
# source sheet this was copied from.
$wsQ3.Range("B2:G4").NumberFormat = "@"

# The copied header used "基金金额"; the new sheet calls this column "基金规模".
$wsQ3.Cells.Item(1,4).Value = "基金规模"

$wsQ3.Cells.Item(2,1).Value = 0
$wsQ3.Cells.Item(2,2).Value = "006165"
$wsQ3.Cells.Item(2,3).Value = "建信中证1000指数增强A"
$wsQ3.Cells.Item(2,4).Value = "3.87"
$wsQ3.Cells.Item(2,5).Value = "84.02"
$wsQ3.Cells.Item(2,6).Value = "1.37"
$wsQ3.Cells.Item(2,7).Value = "0.0530"
$wsQ3.Cells.Item(2,8).Value = 5

$wsQ3.Cells.Item(3,1).Value = 1
$wsQ3.Cells.Item(3,2).Value = "006166"
$wsQ3.Cells.Item(3,3).Value = "建信中证1000指数增强C"
$wsQ3.Cells.Item(3,4).Value = "1.89"
$wsQ3.Cells.Item(3,5).Value = "84.02"
$wsQ3.Cells.Item(3,6).Value = "1.37"
$wsQ3.Cells.Item(3,7).Value = "0.0259"
$wsQ3.Cells.Item(3,8).Value = 5

$wsQ3.Cells.Item(4,1).Value = 2
$wsQ3.Cells.Item(4,2).Value = "013442"
$wsQ3.Cells.Item(4,3).Value = "建信中证1000指数增强E"
$wsQ3.Cells.Item(4,4).Value = "0.18"
$wsQ3.Cells.Item(4,5).Value = "84.02"
$wsQ3.Cells.Item(4,6).Value = "1.37"
$wsQ3.Cells.Item(4,7).Value = "0.0025"
$wsQ3.Cells.Item(4,8).Value = 5

# Drop back to the workbook's default "Normal" style for the data cells now
# that the values are locked in as text -- this matches the unstyled data
# cells used throughout the other per-quarter sheets.
$wsQ3.Range("B2:G4").Style = "Normal"

# Re-apply the index-column style to the new A3/A4 cells (to match A2).
$wsQ3.Range("A2").Copy()
$wsQ3.Range("A3:A4").PasteSpecial(-4122)
$wsQ3.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Restore the previously-active "2021-Q1" sheet as the selected tab, since
#    adding/copying sheets above moved the active tab to the new sheet.
# ---------------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Item(4)
$wsQ1.Activate()

Write-Host "applied 2022-Q3 update"
